$d = $word.ActiveDocument

# Step 1: rename the placeholder paragraph "k" -> "Corruption" (scoped to the
# last paragraph only, with whole-word matching, so we do not touch any other
# "k" characters elsewhere in the document).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Find.Execute("k", $true, $true, $false, $false, $false, $true, 1, $false, "Corruption", 2)

# Step 2: append the new "Corruption" sub-bullets as raw OOXML so the exact
# run/break layout (incl. xml:space="preserve" handling) matches the source.
$endRange = $d.Content
$endRange.Collapse(0)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Once rarely (like, 1/1000) something “corrupts”. This could be:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Upon picking your “egg” up from “daycare”, instead of the normal banter, the music stops. The “daycare worker” says:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>I’m sorry.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The egg we found is…</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Something is wrong. It’s corrupted.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>I can’t belive I’m asking, but you want us to dispose of it, right?</w:t></w:r><w:r><w:br/><w:t xml:space="preserve"> &gt; No | &gt; Yes</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>While trying to load your game from the main screen like normal, instead the music stops and the game freezes. Eventually, the prompt reads:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>I’m terribly sorry.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Your data was corrupted.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>You want to delete the save file, right?</w:t></w:r><w:r><w:br/><w:t>&gt; No | &gt; Yes</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>When selecting “Yes”, the game restarts and everything is normal. Your data is not gone.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>I’m not sure what happens when something is Corrupted, but it should be cool.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Inspired by Terraria’s “Corruption” remix by Qumu (don’t know why; I’ve never played Terraria)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($xml)
